$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab date and update the header string
$ws.Name = "Through 2022-07-06"
$ws.Range("I1").Value = "2022 (through 07-06)"

# Update August (row 8) total for 2022 column
$ws.Range("I8").Value = 30

# Update the grand Total row (row 14) for 2022 column
$ws.Range("I14").Value = 836
